$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at E so the existing "fantasy points" column (E)
# shifts to G, making room for the new "height" / "weight" columns.
$ws.Range("E1:F1").EntireColumn.Insert()

# New header cells, matching the style of the other header cells (B1:D1).
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("E1:F1").Style = $ws.Range("D1").Style

# New data values for each row.
$heights = @(6.5, 6.5, 6.5, 6.5, 6.5, 6.5, 6.5, 6.5, 6.5, 6.5)
$weights = @(251, 251, 251, 251, 251, 251, 251, 251, 251, 251)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $heights[$i]
    $ws.Cells.Item($row, 6).Value = $weights[$i]
}
